$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    "D2" = 18723
    "E2" = 937
    "F2" = 937
    "G2" = 379
    "H2" = 213
    "I2" = 211
    "J2" = 1
    "K2" = 34224
    "L2" = 20944
    "M2" = 13279
    "N2" = 13275
    "O2" = 4
    "P2" = 3688
    "Q2" = 2473
    "R2" = -678
    "S2" = -1844
    "T2" = 985
    "U2" = 1488
    "V2" = 11740
    "W2" = 5
    "X2" = 1.14
    "Y2" = 1.56
    "Z2" = 0.61
    "AA2" = 157.72
    "AB2" = 285.43
    "AC2" = 297
    "AD2" = 78.36
    "AE2" = 19024
    "AF2" = 1.22
    "AG2" = 1000
    "AH2" = 4.3
    "AI2" = 330.27
    "AJ2" = 70133611
    "D3" = 19075
    "E3" = 1340
    "F3" = 1340
    "G3" = 798
    "H3" = 534
    "I3" = 533
    "J3" = 1
    "K3" = 34605
    "L3" = 21359
    "M3" = 13246
    "N3" = 13241
    "O3" = 5
    "P3" = 3688
    "Q3" = 3040
    "R3" = -781
    "S3" = -1234
    "T3" = 1177
    "U3" = 1863
    "V3" = 11220
    "W3" = 7.02
    "X3" = 2.8
    "Y3" = 4.02
    "Z3" = 1.55
    "AA3" = 161.25
    "AB3" = 283.43
    "AC3" = 747
    "AD3" = 31.31
    "AE3" = 18975
    "AF3" = 1.23
    "AG3" = 1000
    "AH3" = 4.27
    "AI3" = 131.1
    "AJ3" = 70133611
    "D4" = 18902
    "E4" = 1240
    "F4" = 1240
    "G4" = 612
    "H4" = 384
    "I4" = 384
    "J4" = 1
    "K4" = 34011
    "L4" = 21092
    "M4" = 12919
    "N4" = 12912
    "O4" = 7
    "P4" = 3688
    "Q4" = 1847
    "R4" = -859
    "S4" = -1024
    "T4" = 1239
    "U4" = 608
    "V4" = 10899
    "W4" = 6.56
    "X4" = 2.03
    "Y4" = 2.94
    "Z4" = 1.12
    "AA4" = 163.26
    "AB4" = 271.98
    "AC4" = 538
    "AD4" = 39.18
    "AE4" = 18504
    "AF4" = 1.14
    "AG4" = 900
    "AH4" = 4.27
    "AI4" = 163.8
    "AJ4" = 70133611
    "D5" = 18899
    "E5" = 872
    "F5" = 872
    "G5" = 397
    "H5" = 127
    "I5" = 127
    "J5" = 1
    "K5" = 35076
    "L5" = 22901
    "M5" = 12175
    "N5" = 12168
    "O5" = 7
    "P5" = 3688
    "Q5" = 3449
    "R5" = -1042
    "S5" = -1058
    "T5" = 923
    "U5" = 2526
    "V5" = 10475
    "W5" = 4.62
    "X5" = 0.67
    "Y5" = 1.01
    "Z5" = 0.37
    "AA5" = 188.1
    "AB5" = 254.67
    "AC5" = 178
    "AD5" = 135.43
    "AE5" = 17437
    "AF5" = 1.38
    "AG5" = 800
    "AH5" = 3.32
    "AI5" = 440.64
    "AJ5" = 70133611
    "D6" = 18856
    "E6" = 904
    "F6" = 904
    "G6" = 420
    "H6" = 223
    "I6" = 222
    "K6" = 34281
    "L6" = 22658
    "M6" = 11624
    "N6" = 11618
    "P6" = 3688
    "Q6" = 1410
    "R6" = -1014
    "S6" = -387
    "T6" = 1211
    "U6" = 199
    "V6" = 10652
    "W6" = 4.8
    "X6" = 1.18
    "Y6" = 1.87
    "Z6" = 0.64
    "AA6" = 194.92
    "AB6" = 240.09
    "AC6" = 312
    "AD6" = 53.17
    "AE6" = 16648
    "AF6" = 1
    "AG6" = 800
    "AH6" = 4.82
    "AI6" = 251.17
    "AJ6" = 70133611
    "D7" = 19930
    "E7" = 945
    "G7" = 250
    "H7" = 56
    "I7" = 76
    "K7" = 34807
    "L7" = 23684
    "M7" = 11123
    "N7" = 11124
    "P7" = 3689
    "Q7" = 1335
    "R7" = -930
    "S7" = -819
    "T7" = 1155
    "U7" = -13
    "W7" = 4.74
    "X7" = 0.28
    "Y7" = 0.67
    "Z7" = 0.16
    "AA7" = 212.92
    "AC7" = 107
    "AD7" = 280.47
    "AE7" = 15941
    "AF7" = 1.88
    "AG7" = 800
    "AH7" = 2.67
    "AI7" = 736.01
    "D8" = 21252
    "E8" = 1586
    "G8" = 1088
    "H8" = 769
    "I8" = 779
    "K8" = 35155
    "L8" = 23850
    "M8" = 11305
    "N8" = 11292
    "P8" = 3689
    "Q8" = 2327
    "R8" = -1080
    "S8" = -734
    "T8" = 1096
    "U8" = 1115
    "W8" = 7.46
    "X8" = 3.62
    "Y8" = 6.95
    "Z8" = 2.2
    "AA8" = 210.97
    "AC8" = 1093
    "AD8" = 27.44
    "AE8" = 16182
    "AF8" = 1.85
    "AG8" = 808
    "AH8" = 2.69
    "AI8" = 72.75
    "D9" = 22231
    "E9" = 1810
    "G9" = 1315
    "H9" = 937
    "I9" = 946
    "K9" = 35861
    "L9" = 24227
    "M9" = 11635
    "N9" = 11600
    "P9" = 3689
    "Q9" = 2333
    "R9" = -1065
    "S9" = -737
    "T9" = 1082
    "U9" = 1262
    "W9" = 8.140000000000001
    "X9" = 4.22
    "Y9" = 8.27
    "Z9" = 2.64
    "AA9" = 208.23
    "AC9" = 1328
    "AD9" = 22.6
    "AE9" = 16623
    "AF9" = 1.8
    "AG9" = 817
    "AH9" = 2.72
    "AI9" = 60.54
}

foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}
